# Fruta / hortaliza, semanal
# Insert 3 new weekly rows of "Palta" price data (Hass, Primera/Segunda/Tercera from Peru)
# at the top of this product's data block (rows 203-205), pushing the existing
# rows 203-229 down to 206-232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current row 203 -- this shifts the
# existing data (rows 203:229) down to rows 206:232, exactly like Excel's
# normal "Insert Rows" behaviour (formats/styles copy down from the row below).
$ws.Rows("203:205").Insert()

# Helper to populate one new data row. Columns A,B,C,E,F,G,H,I,J are the
# constant "boilerplate" for this sheet (same market/product on every row).
# NOTE: uses positional parameters (named "-Foo bar" args are not reliable
# for user-defined functions in this COM-interop runtime).
function Set-PaltaRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {

    $ws.Cells.Item($Row, 1).Value2 = 1
    $ws.Cells.Item($Row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($Row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($Row, 4).Value2 = $Fecha
    $ws.Cells.Item($Row, 5).Value2 = 15
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value2 = 100106
    $ws.Cells.Item($Row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($Row, 9).Value2 = 100106002
    $ws.Cells.Item($Row, 10).Value = "Palta"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value2 = $Volumen
    $ws.Cells.Item($Row, 14).Value2 = $PrecioMin
    $ws.Cells.Item($Row, 15).Value2 = $PrecioMax
    $ws.Cells.Item($Row, 16).Value2 = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value2 = $PrecioKg
    $ws.Cells.Item($Row, 20).Value2 = $KgUnidad
}

Set-PaltaRow 203 45124 "Hass" "Primera" 200 23000 24000 23500 "$/bandeja 10 kilos" "Perú" 2350 10
Set-PaltaRow 204 45124 "Hass" "Segunda" 250 21000 22000 21600 "$/bandeja 10 kilos" "Perú" 2160 10
Set-PaltaRow 205 45124 "Hass" "Tercera" 250 19000 20000 19600 "$/bandeja 10 kilos" "Perú" 1960 10

# Apply the same date number format used by the other "Fecha" cells (style index 2 / numFmtId 165)
$ws.Range("D203:D205").NumberFormat = $ws.Range("D206").NumberFormat
